$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 97, pushing the old rows 97-127 down to 99-129.
$ws.Range("A97:A98").EntireRow.Insert()

# Remember the date number format used by the date column so the two new rows
# keep the same formatting as every other data row.
$dateFmt = $ws.Range("D99").NumberFormat

# Populate the two new rows with the new weekly data. Columns that are constant
# across every data row in this sheet (A, B, C, E, F, G, H, I, J, K, Q, R, T) are
# copied verbatim from the neighboring rows.
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = 'Femacal de La Calera'
$ws.Range("C97").Value = 'Coquimbo'
$ws.Range("D97").Value = 44489
$ws.Range("D97").NumberFormat = $dateFmt
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = 'Fruta'
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = 'Berries'
$ws.Range("I97").Value = 100112025
$ws.Range("J97").Value = 'Frutilla'
$ws.Range("K97").Value = 'Sin especificar'
$ws.Range("L97").Value = 'Especial'
$ws.Range("M97").Value = 97
$ws.Range("N97").Value = 9000
$ws.Range("O97").Value = 10000
$ws.Range("P97").Value = 9515
$ws.Range("Q97").Value = '$/bandeja 7 kilos'
$ws.Range("R97").Value = 'Provincia de Melipilla'
$ws.Range("S97").Value = 1359
$ws.Range("T97").Value = 7

$ws.Range("A98").Value = 3
$ws.Range("B98").Value = 'Femacal de La Calera'
$ws.Range("C98").Value = 'Coquimbo'
$ws.Range("D98").Value = 44489
$ws.Range("D98").NumberFormat = $dateFmt
$ws.Range("E98").Value = 5
$ws.Range("F98").Value = 'Fruta'
$ws.Range("G98").Value = 100101
$ws.Range("H98").Value = 'Berries'
$ws.Range("I98").Value = 100112025
$ws.Range("J98").Value = 'Frutilla'
$ws.Range("K98").Value = 'Sin especificar'
$ws.Range("L98").Value = 'Segunda'
$ws.Range("M98").Value = 50
$ws.Range("N98").Value = 7000
$ws.Range("O98").Value = 7000
$ws.Range("P98").Value = 7000
$ws.Range("Q98").Value = '$/bandeja 7 kilos'
$ws.Range("R98").Value = 'Provincia de Melipilla'
$ws.Range("S98").Value = 1000
$ws.Range("T98").Value = 7
